# Task: Completed daily operations, 8 hours, 10/06
# Adds a new Time Log entry row (row 11) for 10/06/2023:
#   Date = 10/06/2023 (serial 45205), Name of Task = Internship,
#   Description = "Contributed technical work by aiding in resolving
#   inconsistencies flagged by the system for employee calls"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 11

# Date column - match the date formatting/style already used by the
# previous entries (column A) by copying the number format from the
# last existing data row.
$ws.Cells.Item($newRow, 1).Value = 45205
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat

# Name of Task column
$ws.Cells.Item($newRow, 2).Value = "Internship"

# Description column
$ws.Cells.Item($newRow, 3).Value = "Contributed technical work by aiding in resolving inconsistencies flagged by the system for employee calls"

# Move the active selection below the newly entered row, matching the
# cursor position left behind after typing the row and pressing Enter.
$ws.Cells.Item($newRow + 1, 3).Select() | Out-Null
